$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '301.64'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.83%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '32.86'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '4.21%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.957'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-2.71%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07774'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-1.09%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.961'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-13.88%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.848'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.53%'
$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.800'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-1.34%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9232'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.18%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1769'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '1.28%'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07882'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '4.07%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08647'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-7.10%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03150'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '4.69%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.1002'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.02%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001515'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.41%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005932'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.74%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.461'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.34%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.154'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-3.96%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3338'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2.03%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1319'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.70%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.308'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '9.46%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '16.38%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04564'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-1.23%'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-2.35%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004436'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.76%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.03%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01714'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-1.38%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04726'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2.42%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007679'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '10.01%'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-0.22%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002343'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '6.88%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01056'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '8.49%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00006265'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-0.40%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.04%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.8206'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '10.15%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003104'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-61.13%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002103'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.04%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.04%'
